# finish(change) the ImportExcel and Change the importExcel
# - rename sheets: Sheet1->gg, Sheet2->sc, Sheet3->qq
# - fill in sample data on "sc" (formerly Sheet2)
# - move the active tab/selection onto "qq" (formerly Sheet3)

$wb = $excel.ActiveWorkbook

# --- rename the three sheets ---
$wb.Worksheets.Item(1).Name = "gg"
$wb.Worksheets.Item(2).Name = "sc"
$wb.Worksheets.Item(3).Name = "qq"

$wsGG = $wb.Worksheets.Item("gg")
$wsSC = $wb.Worksheets.Item("sc")
$wsQQ = $wb.Worksheets.Item("qq")

# --- keep "gg" selection where it was (B6), it just stops being the active tab ---
[void]$wsGG.Range("B6").Select()

# --- populate "sc" with the new table data ---
$wsSC.Range("A1").Value = "c"
$wsSC.Range("B1").Value = "s"
$wsSC.Range("A2").Value = "sdf"
$wsSC.Range("B2").Value = "sdf"
$wsSC.Range("A3").Value = "sdf"
$wsSC.Range("B3").Value = "sdf"
[void]$wsSC.Range("B3").Select()

# --- move the live selection to "qq" and make it the active tab ---
[void]$wsQQ.Range("E37").Select()
$wsQQ.Activate()
